$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing literal text storage. The source
# workbook keeps every Price/Volume cell as inline text (even values that
# look numeric, e.g. "338.75" or "30.192.19"), so a plain .Value assignment
# would let Excel auto-coerce them to numbers (dropping trailing zeros like
# "0.09330" -> 0.0933, and changing the stored cell type). Briefly tagging
# the cell as Text (@) before the write avoids that coercion, and the
# trailing ClearFormats() removes the temporary format override again so no
# stray style survives in the saved workbook.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue 'D2' '30.192.19'
$ws.Range('E2').Value = '  +0.79%  '
Set-TextValue 'D3' '2.080.05'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('E4').Value = '  -0.64%  '
Set-TextValue 'D5' '338.75'
$ws.Range('E6').Value = '  -0.49%  '
Set-TextValue 'D7' '0.5258'
$ws.Range('E7').Value = '  +1.04%  '
Set-TextValue 'D8' '0.4357'
$ws.Range('E8').Value = '  -2.10%  '
Set-TextValue 'D9' '54.89'
$ws.Range('E9').Value = '  +0.76%  '
Set-TextValue 'D10' '0.09330'
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('E11').Value = '  -0.69%  '
Set-TextValue 'D12' '24.47'
$ws.Range('E12').Value = '  -2.73%  '
Set-TextValue 'D13' '8.457'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '2.097.01'
$ws.Range('E14').Value = '  -2.86%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D15' '6.848'
$ws.Range('E15').Value = '  -0.38%  '
Set-TextValue 'D16' '100.24'
$ws.Range('E16').Value = '  -1.99%  '
Set-TextValue 'D17' '0.00001157'
$ws.Range('E17').Value = '  -0.53%  '
Set-TextValue 'D18' '1.004'
$ws.Range('E18').Value = '  -0.52%  '
Set-TextValue 'D19' '20.84'
$ws.Range('E19').Value = '  -3.23%  '
Set-TextValue 'D20' '0.06700'
$ws.Range('E20').Value = '  +0.16%  '
Set-TextValue 'D21' '6.294'
$ws.Range('E21').Value = '  +0.18%  '
Set-TextValue 'D22' '1.002'
$ws.Range('E22').Value = '  -0.45%  '
Set-TextValue 'D23' '30.208.40'
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('E24').Value = '  -2.87%  '
Set-TextValue 'D25' '2.312'
$ws.Range('E25').Value = '  -0.79%  '
Set-TextValue 'D26' '21.72'
$ws.Range('E26').Value = '  -1.72%  '
Set-TextValue 'D27' '162.34'
$ws.Range('E27').Value = '  -0.19%  '
Set-TextValue 'D28' '6.780'
$ws.Range('E28').Value = '  +3.15%  '
Set-TextValue 'D29' '2.479'
$ws.Range('E29').Value = '  -2.69%  '
Set-TextValue 'D30' '133.13'
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('E31').Value = '  -2.15%  '
$ws.Range('E32').Value = '  -6.76%  '
Set-TextValue 'D33' '0.1045'
$ws.Range('E33').Value = '  -1.06%  '
Set-TextValue 'D34' '6.237'
$ws.Range('E34').Value = '  -0.05%  '
Set-TextValue 'D35' '3.908'
$ws.Range('E35').Value = '  -1.57%  '
Set-TextValue 'D36' '0.02601'
$ws.Range('E36').Value = '  -0.09%  '
Set-TextValue 'D37' '9.831'
$ws.Range('E37').Value = '  -8.99%  '
Set-TextValue 'D38' '0.06706'
$ws.Range('E38').Value = '  -2.19%  '
Set-TextValue 'D39' '0.6940'
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('E40').Value = '  -1.57%  '
Set-TextValue 'D41' '1.326'
$ws.Range('E41').Value = '  -0.62%  '
Set-TextValue 'D42' '0.2197'
$ws.Range('E42').Value = '  -2.29%  '
Set-TextValue 'D43' '0.6706'
$ws.Range('E43').Value = '  -1.99%  '
Set-TextValue 'D44' '2.357'
$ws.Range('E44').Value = '  +0.23%  '
Set-TextValue 'D45' '14.22'
$ws.Range('E45').Value = '  -1.93%  '
$ws.Range('E46').Value = '  -0.43%  '
Set-TextValue 'D47' '1.310'
$ws.Range('E47').Value = '  +4.98%  '
Set-TextValue 'D48' '3.627'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('E49').Value = '  -2.31%  '
Set-TextValue 'D50' '1.209'
$ws.Range('E50').Value = '  +2.01%  '
Set-TextValue 'D51' '1.205'
$ws.Range('E51').Value = '  -1.47%  '
